# Edit script: insert two new data columns (D,E) into the SYBT quarterly financials sheet,
# shifting the existing quarters right, and populate the two new columns with the newest
# quarter's figures (commit message: "Doing Updates for Financials").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SYBT")

# 1) Insert two blank columns at D:E (existing D:K shifts right to F:M).
$ws.Range("D7:E102").Insert(-4161)

# 2) Copy number formats from the (now shifted) old "D/E" columns, F:G, into the
#    new D:E columns so the new cells pick up the same date / number styling.
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rowData = @{
    7 = @(43465, 43373)
    8 = @(35000, 33000)
    9 = @("NA", "NA")
    10 = @("NA", "NA")
    11 = @($null, $null)
    12 = @("NA", "NA")
    13 = @(0, 0)
    14 = @(0, 0)
    15 = @(0, 0)
    16 = @($null, $null)
    17 = @(5100, 5200)
    18 = @(29900, 27800)
    19 = @($null, $null)
    20 = @(-13000, -10300)
    21 = @(18700, 18700)
    22 = @(0, 0)
    23 = @(16900, 17400)
    24 = @(2300, 3600)
    25 = @(0, 0)
    26 = @(14700, 13900)
    27 = @(14700, 13900)
    28 = @(0, 0)
    29 = @("NA", "NA")
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(13000, 10300)
    33 = @(14700, 13900)
    34 = @(0, 0)
    35 = @(14700, 13900)
    38 = @(43465, 43373)
    39 = @($null, $null)
    40 = @($null, $null)
    41 = @(51900, 66000)
    42 = @(157400, 64800)
    43 = @(0, 0)
    44 = @(0, 0)
    45 = @(0, 0)
    46 = @(0, 0)
    47 = @(0, 0)
    48 = @(44800, 43600)
    49 = @(0, 0)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(0, 0)
    53 = @(0, 0)
    54 = @(3302900, 3324800)
    55 = @($null, $null)
    56 = @($null, $null)
    57 = @(800, 700)
    58 = @(0, 0)
    59 = @(0, 0)
    60 = @(0, 0)
    61 = @(0, 0)
    62 = @(0, 0)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(2936400, 2971800)
    67 = @($null, $null)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(298200, 289300)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(366500, 353000)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(14700, 13900)
    82 = @($null, $null)
    83 = @(1800, 1300)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(20200, 12200)
    90 = @($null, $null)
    91 = @(-2100, -2200)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(106900, 61800)
    95 = @($null, $null)
    96 = @(-5700, -5700)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(-48600, -8500)
    101 = @(0, 0)
    102 = @(78500, 65500)
}

foreach ($r in $rowData.Keys) {
    $pair = $rowData[$r]
    $dVal = $pair[0]
    $eVal = $pair[1]
    if ($null -ne $dVal) {
        $ws.Cells.Item([int]$r, 4).Value = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item([int]$r, 5).Value = $eVal
    }
}
